$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.869.88'
$ws.Range("E2").Value = '  -0.52%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.670.30'
$ws.Range("E3").Value = '  +1.07%  '

$ws.Range("E4").Value = '  +0.40%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.76'
$ws.Range("E5").Value = '  +0.10%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.531'
$ws.Range("E6").Value = '  +4.91%  '

$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("E8").Value = '  +2.50%  '

$ws.Range("E9").Value = '  +0.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.31'
$ws.Range("E10").Value = '  +3.93%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0894'
$ws.Range("E11").Value = '  +4.10%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.907.32'
$ws.Range("E12").Value = '  +1.15%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.694.30'
$ws.Range("E13").Value = '  +2.51%  '

$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.64'
$ws.Range("E16").Value = '  +1.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.914.86'
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.52'
$ws.Range("E18").Value = '  -4.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.81'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0734'
$ws.Range("E20").Value = '  +0.68%  '

$ws.Range("E21").Value = '  +0.29%  '

$ws.Range("E22").Value = '  +1.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.20'
$ws.Range("E23").Value = '  -0.68%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  -2.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.47'
$ws.Range("E25").Value = '  -0.57%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.15'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.116'
$ws.Range("E27").Value = '  +2.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.95'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("E29").Value = '  +0.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("E31").Value = '  +0.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +1.51%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.464.83'
$ws.Range("E33").Value = '  -3.32%  '

$ws.Range("E34").Value = '  +3.70%  '

$ws.Range("E35").Value = '  +3.08%  '

$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.899'
$ws.Range("E37").Value = '  +2.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.570'
$ws.Range("E38").Value = '  -1.32%  '

$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.89'
$ws.Range("E40").Value = '  -1.98%  '

$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("E42").Value = '  +3.07%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.78'
$ws.Range("E43").Value = '  +1.23%  '

$ws.Range("E44").Value = '  +6.48%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.816.52'
$ws.Range("E45").Value = '  +1.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.780'
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.53'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  +9.34%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.100'
$ws.Range("E50").Value = '  +2.41%  '

$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0508'
$ws.Range("E51").Value = '  +0.99%  '
